# Implement first version of lot sizing rules:
# - Increase NrBuckets (Generic sheet, B4) from 4 to 5
# - Add a 6th row (bucket "4") to ForecastedAverageDemand and
#   ForcastedStandardDeviation sheets, repeating the existing bucket pattern.

$wb = $excel.ActiveWorkbook

# --- Generic sheet: bump NrBuckets from 4 to 5 ---
$wsGeneric = $wb.Worksheets.Item("Generic")
$wsGeneric.Range("B4").Value = 5

# --- ForecastedAverageDemand sheet: append row 6 ---
$wsAvg = $wb.Worksheets.Item("ForecastedAverageDemand")
# Carry the bold/bordered "bucket index" formatting from A5 down to A6.
$wsAvg.Range("A5").Copy($wsAvg.Range("A6"))
$wsAvg.Cells.Item(6, 1).Value = 4
$wsAvg.Cells.Item(6, 2).Value = 0
$wsAvg.Cells.Item(6, 3).Value = 600
$wsAvg.Cells.Item(6, 4).Value = 150
$wsAvg.Cells.Item(6, 5).Value = 0
$wsAvg.Cells.Item(6, 6).Value = 300
$wsAvg.Cells.Item(6, 7).Value = 90
$wsAvg.Cells.Item(6, 8).Value = 30
$wsAvg.Cells.Item(6, 9).Value = 420
$wsAvg.Cells.Item(6, 10).Value = 300
$wsAvg.Cells.Item(6, 11).Value = 900
$wsAvg.Cells.Item(6, 12).Value = 210
$wsAvg.Cells.Item(6, 13).Value = 0
$wsAvg.Cells.Item(6, 14).Value = 0
$wsAvg.Cells.Item(6, 15).Value = 0
$wsAvg.Cells.Item(6, 16).Value = 0
$wsAvg.Cells.Item(6, 17).Value = 0
$wsAvg.Cells.Item(6, 18).Value = 0
$wsAvg.Cells.Item(6, 19).Value = 0
$wsAvg.Cells.Item(6, 20).Value = 0
$wsAvg.Cells.Item(6, 21).Value = 0
$wsAvg.Cells.Item(6, 22).Value = 0
$wsAvg.Cells.Item(6, 23).Value = 0

# --- ForcastedStandardDeviation sheet: append row 6 ---
$wsStd = $wb.Worksheets.Item("ForcastedStandardDeviation")
# Carry the bold/bordered "bucket index" formatting from A5 down to A6.
$wsStd.Range("A5").Copy($wsStd.Range("A6"))
$wsStd.Cells.Item(6, 1).Value = 4
$wsStd.Cells.Item(6, 2).Value = 0
$wsStd.Cells.Item(6, 3).Value = 5
$wsStd.Cells.Item(6, 4).Value = 5
$wsStd.Cells.Item(6, 5).Value = 0
$wsStd.Cells.Item(6, 6).Value = 5
$wsStd.Cells.Item(6, 7).Value = 5
$wsStd.Cells.Item(6, 8).Value = 5
$wsStd.Cells.Item(6, 9).Value = 2
$wsStd.Cells.Item(6, 10).Value = 3
$wsStd.Cells.Item(6, 11).Value = 5
$wsStd.Cells.Item(6, 12).Value = 5
$wsStd.Cells.Item(6, 13).Value = 0
$wsStd.Cells.Item(6, 14).Value = 0
$wsStd.Cells.Item(6, 15).Value = 0
$wsStd.Cells.Item(6, 16).Value = 0
$wsStd.Cells.Item(6, 17).Value = 0
$wsStd.Cells.Item(6, 18).Value = 0
$wsStd.Cells.Item(6, 19).Value = 0
$wsStd.Cells.Item(6, 20).Value = 0
$wsStd.Cells.Item(6, 21).Value = 0
$wsStd.Cells.Item(6, 22).Value = 0
$wsStd.Cells.Item(6, 23).Value = 0
